$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "41.602.33"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.72%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.169.87"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -2.95%  "
$ws.Range("E4").Value = "  -0.21%  "
$ws.Range("E5").Value = "  -1.99%  "
$ws.Range("E6").Value = "  -3.15%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "72.19"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -3.41%  "
$ws.Range("E8").Value = "  -0.17%  "
$ws.Range("E9").Value = "  -5.22%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "40.02"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -6.77%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0908"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -5.64%  "
$ws.Range("E12").Value = "  -4.26%  "
$ws.Range("E13").Value = "  -3.41%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.72"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -4.29%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.492.33"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.28%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.27"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.62%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.159.31"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.70%  "
$ws.Range("E18").Value = "  -7.07%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "41.468.50"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.82%  "
$ws.Range("E20").Value = "  -3.02%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "70.09"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -4.22%  "
$ws.Range("E22").Value = "  -7.33%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.78"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -13.83%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "226.75"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.02%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.01"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -4.22%  "
$ws.Range("E26").Value = "  +0.16%  "
$ws.Range("E27").Value = "  -6.42%  "
$ws.Range("E28").Value = "  -9.83%  "
$ws.Range("E29").Value = "  -4.15%  "
$ws.Range("E30").Value = "  -1.42%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "171.15"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.54%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "19.81"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.99%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "33.33"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +8.66%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0773"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.79%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.23"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -7.95%  "
$ws.Range("E36").Value = "  -3.71%  "
$ws.Range("E37").Value = "  -1.79%  "
$ws.Range("E38").Value = "  -4.95%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0308"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.79%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "12.16"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -8.26%  "
$ws.Range("E41").Value = "  -2.20%  "
$ws.Range("E42").Value = "  -5.98%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "58.90"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -9.41%  "
$ws.Range("E44").Value = "  -3.51%  "
$ws.Range("E45").Value = "  -5.48%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0967"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -4.08%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "97.81"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -6.79%  "
$ws.Range("E48").Value = "  -3.90%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.19"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -7.61%  "
